# Updated cryptos list on Wed Sep 27 19:16:46 UTC 2023 with GitHub Actions
#
# NOTE: several "Price" (column D) values are plain decimal-looking strings
# (e.g. "211.72"). Excel's COM automation auto-detects such strings and
# stores them as numbers, but the source workbook stores every Price cell
# as text. To faithfully reproduce a text cell (just like typing '211.72
# into Excel) we prefix those values with a leading apostrophe, which
# forces text storage while keeping the displayed/read value identical.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 19 / Row 20 swap (ShibaInu <-> Chainlink) ---
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "'7.65"
$ws.Range("E19").Value = "  +5.27%  "

$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.0₃0721"
$ws.Range("E20").Value = "  -0.23%  "

# --- Row 23 / Row 24 swap (Avalanche <-> Toncoin) ---
$ws.Range("B23").Value = "Toncoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D23").Value = "'2.14"
$ws.Range("E23").Value = "  +1.15%  "

$ws.Range("B24").Value = "Avalanche"
$ws.Range("C24").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D24").Value = "'8.92"
$ws.Range("E24").Value = "  -0.52%  "

# --- Price (D) and Volume(1h) (E) updates for all other rows ---

$ws.Range("D2").Value = "26.262.62"
$ws.Range("E2").Value = "  +0.28%  "

$ws.Range("D3").Value = "1.594.40"
$ws.Range("E3").Value = "  +0.57%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "'211.72"
$ws.Range("E5").Value = "  +0.09%  "

$ws.Range("D6").Value = "'0.505"
$ws.Range("E6").Value = "  -0.08%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("E9").Value = "  +0.41%  "

$ws.Range("D10").Value = "'18.96"
$ws.Range("E10").Value = "  -1.08%  "

$ws.Range("D11").Value = "'0.0852"
$ws.Range("E11").Value = "  +0.63%  "

$ws.Range("D12").Value = "1.818.69"

$ws.Range("D13").Value = "1.593.32"
$ws.Range("E13").Value = "  +0.27%  "

$ws.Range("E14").Value = "  -0.10%  "

$ws.Range("E15").Value = "  -2.42%  "

$ws.Range("D16").Value = "'63.62"
$ws.Range("E16").Value = "  -0.37%  "

$ws.Range("D17").Value = "26.229.17"
$ws.Range("E17").Value = "  +0.20%  "

$ws.Range("D18").Value = "'229.76"
$ws.Range("E18").Value = "  +7.43%  "

$ws.Range("E21").Value = "  +0.06%  "

$ws.Range("E22").Value = "  -0.19%  "

$ws.Range("D25").Value = "'145.62"
$ws.Range("E25").Value = "  +0.99%  "

$ws.Range("E26").Value = "  -0.01%  "

$ws.Range("E27").Value = "  +0.25%  "

$ws.Range("E28").Value = "  +0.75%  "

$ws.Range("D29").Value = "'15.32"
$ws.Range("E29").Value = "  +1.63%  "

$ws.Range("E31").Value = "  +0.15%  "

$ws.Range("E32").Value = "  +0.74%  "

$ws.Range("D33").Value = "1.463.13"
$ws.Range("E33").Value = "  +3.44%  "

$ws.Range("E34").Value = "  +0.45%  "

$ws.Range("E35").Value = "  -0.54%  "

$ws.Range("E36").Value = "  +0.64%  "

$ws.Range("D37").Value = "'0.566"
$ws.Range("E37").Value = "  -3.32%  "

$ws.Range("E38").Value = "  -1.13%  "

$ws.Range("E39").Value = "  +0.14%  "

$ws.Range("E40").Value = "  -1.83%  "

$ws.Range("E41").Value = "  +0.02%  "

$ws.Range("E42").Value = "  +2.12%  "

$ws.Range("E43").Value = "  -1.52%  "

$ws.Range("D44").Value = "1.731.46"
$ws.Range("E44").Value = "  +0.71%  "

$ws.Range("D45").Value = "'0.757"
$ws.Range("E45").Value = "  -0.96%  "

$ws.Range("D46").Value = "'60.44"
$ws.Range("E46").Value = "  -0.64%  "

$ws.Range("D47").Value = "'87.58"
$ws.Range("E47").Value = "  +2.74%  "

$ws.Range("E48").Value = "  -0.42%  "

$ws.Range("E49").Value = "  +0.04%  "

$ws.Range("E50").Value = "  -0.02%  "

$ws.Range("D51").Value = "'0.0946"
$ws.Range("E51").Value = "  -2.24%  "
